# Update countries & provincias Spain
#
# 1) Swap the shared-string labels for two pairs of rows (the sorted-by-
#    cases list re-ranked them): "San Marino" / "Curazao" and
#    "Islas Malvinas" / "Montserrat".
# 2) Refresh the case-count columns (B:H) for the rows whose underlying
#    numbers changed, including the two rows above whose data moved one
#    slot down/up together with the label swap.
#
# NOTE: this COM-interop dialect does not bind PowerShell *named*
# parameters (-Foo bar) on user functions, so Set-Row takes its
# arguments positionally; pass $null for any column that should be
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($Row, $Pais, $Total, $Nuevos, $Activos, $Recuperados, $Criticos, $MuertesHoy, $Muertes)

    if ($null -ne $Pais) { $ws.Cells.Item($Row, 1).Value = $Pais }
    if ($null -ne $Total) { $ws.Cells.Item($Row, 2).Value = $Total }
    if ($null -ne $Nuevos) { $ws.Cells.Item($Row, 3).Value = $Nuevos }
    if ($null -ne $Activos) { $ws.Cells.Item($Row, 4).Value = $Activos }
    if ($null -ne $Recuperados) { $ws.Cells.Item($Row, 5).Value = $Recuperados }
    if ($null -ne $Criticos) { $ws.Cells.Item($Row, 6).Value = $Criticos }
    if ($null -ne $MuertesHoy) { $ws.Cells.Item($Row, 7).Value = $MuertesHoy }
    if ($null -ne $Muertes) { $ws.Cells.Item($Row, 8).Value = $Muertes }
}

# Estados Unidos
Set-Row 4 $null 8819311 72028 5734149 2855155 $null 723 230007

# India
Set-Row 5 $null 7863892 50224 7075723 669602 $null 575 118567

# Brasil
Set-Row 6 $null 5380635 24985 $null 425860 $null 375 156903

# Rank 24 (Alemania)
Set-Row 20 $null 427799 10449 $null 103588 $null 21 10111

# Rank 62 (Barein)
Set-Row 58 $null 79975 401 76474 3189 $null $null $null

# Rank 88 (Bulgaria)
Set-Row 84 $null 37562 1043 18232 18246 $null 7 1084

# Rank 109 (Namibia)
Set-Row 105 $null 12579 78 10775 1671 $null $null $null

# Rank 122 (Zimbabue)
Set-Row 118 $null 8269 12 7785 248 $null $null $null

# Rank 124 (Mauritania)
Set-Row 120 $null 7663 1 7378 122 $null $null $null

# Rank 131 (Suazilandia)
Set-Row 127 $null 5847 16 5491 240 $null $null $null

# Rank 137 (Siria)
Set-Row 133 $null 5359 40 1722 3370 $null 3 267

# Rank 158 (Sudan del Sur)
Set-Row 154 $null 2878 2 $null 1532 $null $null $null

# Rank 174 slot: now Curazao (updated numbers, moved up from the 175 slot)
Set-Row 170 "Curazao" 837 19 555 281 $null $null 1

# Rank 175 slot: now San Marino (its own numbers are unchanged, just shifted down)
Set-Row 171 "San Marino" 819 $null 716 61 $null $null 42

# Rank 220 slot: now Montserrat
Set-Row 216 "Montserrat" $null $null 12 $null $null $null 1

# Rank 221 slot: now Islas Malvinas
Set-Row 217 "Islas Malvinas" $null $null 13 $null $null $null 0
